$wb = $excel.ActiveWorkbook

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5748569
$ws.Range("I61").Value = 5953732.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 5953732.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -5953520.5
$ws.Range("N61").Value = -4424
$ws.Range("H63").Value = 2939.4
$ws.Range("I63").Value = 3050.261
$ws.Range("J63").Value = 2575.1428
$ws.Range("K63").Value = 3050.261
$ws.Range("L63").Value = 2575.1428
$ws.Range("M63").Value = -2364.261
$ws.Range("N63").Value = -3947.1428
$ws.Range("H66").Value = 2939.4
$ws.Range("I66").Value = 3050.261
$ws.Range("J66").Value = 2575.1428
$ws.Range("K66").Value = 15251.305
$ws.Range("L66").Value = 12875.714
$ws.Range("M66").Value = -11819.305
$ws.Range("N66").Value = -19739.714
$ws.Range("H88").Value = 2747
$ws.Range("I88").Value = 2799.9092
$ws.Range("J88").Value = 2650
$ws.Range("K88").Value = 2799.9092
$ws.Range("L88").Value = 2650
$ws.Range("M88").Value = -2393.9092
$ws.Range("N88").Value = -3462
$ws.Range("H91").Value = 2747
$ws.Range("I91").Value = 2799.9092
$ws.Range("J91").Value = 2650
$ws.Range("K91").Value = 2799.9092
$ws.Range("L91").Value = 2650
$ws.Range("M91").Value = -1395.9092
$ws.Range("N91").Value = -5458
$ws.Range("H132").Value = 1436016.5
$ws.Range("I132").Value = 1069.5333
$ws.Range("J132").Value = 5349508
$ws.Range("K132").Value = 3208.5999
$ws.Range("L132").Value = 16048524
$ws.Range("M132").Value = -678.5999000000002
$ws.Range("N132").Value = -16053584
$ws.Range("H136").Value = 5748569
$ws.Range("I136").Value = 5953732.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 17861197.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -17858647.5
$ws.Range("N136").Value = -17100

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 26635
$ws.Range("J55").Value = 26635
$ws.Range("L55").Value = 26635
$ws.Range("N55").Value = -27181

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2366.6667
$ws.Range("I15").Value = 1200
$ws.Range("J15").Value = 2950
$ws.Range("K15").Value = 1200
$ws.Range("L15").Value = 2950
$ws.Range("M15").Value = -1030
$ws.Range("N15").Value = -3290

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 4200
$ws.Range("J74").Value = 4200
$ws.Range("L74").Value = 12600
$ws.Range("N74").Value = -14722
$ws.Range("H77").Value = 4200
$ws.Range("J77").Value = 4200
$ws.Range("L77").Value = 37800
$ws.Range("N77").Value = -48408
$ws.Range("H123").Value = 2955
$ws.Range("I123").Value = 2930
$ws.Range("J123").Value = 2980
$ws.Range("K123").Value = 8790
$ws.Range("L123").Value = 8940
$ws.Range("M123").Value = -6340
$ws.Range("N123").Value = -13840
$ws.Range("H125").Value = 4000
$ws.Range("J125").Value = 4000
$ws.Range("L125").Value = 12000
$ws.Range("N125").Value = -21840
$ws.Range("H137").Value = 29413784
$ws.Range("I137").Value = 55556308
$ws.Range("J137").Value = 3445.375
$ws.Range("K137").Value = 166668924
$ws.Range("L137").Value = 10336.125
$ws.Range("M137").Value = -166663824
$ws.Range("N137").Value = -20536.125

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5264944
$ws.Range("I80").Value = 1885.5555
$ws.Range("J80").Value = 100000000
$ws.Range("K80").Value = 1885.5555
$ws.Range("L80").Value = 100000000
$ws.Range("M80").Value = -887.5554999999999
$ws.Range("N80").Value = -100001996
$ws.Range("H83").Value = 5264944
$ws.Range("I83").Value = 1885.5555
$ws.Range("J83").Value = 100000000
$ws.Range("K83").Value = 9427.7775
$ws.Range("L83").Value = 500000000
$ws.Range("M83").Value = -4435.7775
$ws.Range("N83").Value = -500009984
$ws.Range("H92").Value = 6572
$ws.Range("J92").Value = 6572
$ws.Range("L92").Value = 6572
$ws.Range("N92").Value = -10316
$ws.Range("H94").Value = 9448
$ws.Range("J94").Value = 9448
$ws.Range("L94").Value = 9448
$ws.Range("N94").Value = -10800
$ws.Range("H95").Value = 22366.666
$ws.Range("J95").Value = 22366.666
$ws.Range("L95").Value = 22366.666
$ws.Range("N95").Value = -27858.666
$ws.Range("H96").Value = 10754.25
$ws.Range("J96").Value = 10754.25
$ws.Range("L96").Value = 10754.25
$ws.Range("N96").Value = -16246.25
$ws.Range("H97").Value = 1766.3334
$ws.Range("I97").Value = 1766.3334
$ws.Range("K97").Value = 1766.3334
$ws.Range("M97").Value = -1270.3334
$ws.Range("H98").Value = 22482.25
$ws.Range("J98").Value = 22482.25
$ws.Range("L98").Value = 22482.25
$ws.Range("N98").Value = -28472.25
$ws.Range("H99").Value = 7560
$ws.Range("I99").Value = 3433.3333
$ws.Range("J99").Value = 13750
$ws.Range("K99").Value = 3433.3333
$ws.Range("L99").Value = 13750
$ws.Range("M99").Value = -1187.3333
$ws.Range("N99").Value = -18242
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490
$ws.Range("H102").Value = 1442.7826
$ws.Range("I102").Value = 909.2
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 909.2
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 712.8
$ws.Range("N102").Value = -8244
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 30693610
$ws.Range("I55").Value = 2293701.8
$ws.Range("J55").Value = 66666828
$ws.Range("K55").Value = 2293701.8
$ws.Range("L55").Value = 66666828
$ws.Range("M55").Value = -2293528.8
$ws.Range("N55").Value = -66667174

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 680
$ws.Range("I126").Value = 240
$ws.Range("K126").Value = 720
$ws.Range("M126").Value = 1750
